$wb = $excel.ActiveWorkbook

# --- Work on the "emission" sheet (values + formulas in row 2) ---
$ws = $wb.Worksheets.Item("emission")

# C2:G2 become formulas (C2..F2 step halving, G2 = F2*0.7)
$ws.Range("C2").Formula = "=B2*0.5"
$ws.Range("D2").Formula = "=C2*0.5"
$ws.Range("E2").Formula = "=D2*0.5"
$ws.Range("F2").Formula = "=E2*0.5"
$ws.Range("G2").Formula = "=F2*0.7"

# H2:AA2 become flat static values of 200000
$ws.Range("H2:AA2").Value = 200000

# --- Update view / selection state ---
# Make "emission" the active (selected) sheet/tab
$ws.Activate()

# Scroll so column G is the left-most visible column, then select H2:AA2
$excel.ActiveWindow.TopLeftCell = $ws.Range("G1")
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("H2:AA2").Select()
